$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entregables")

# Center the status column (C2:C10 and C12:C18) and set cell C9 to "x"
$range1 = $ws.Range("C2:C10")
$range1.HorizontalAlignment = -4108  # xlCenter
$range1.VerticalAlignment = -4108    # xlCenter

$range2 = $ws.Range("C12:C18")
$range2.HorizontalAlignment = -4108
$range2.VerticalAlignment = -4108

$ws.Range("C9").Value = "x"

$ws.Range("C9").Select()
